$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scenario (B), year (C), value (D) for all data rows 2-25.
# Rows 23-25 are brand new rows; also set A (model), E (region), F (unit), G (variable).
$ws.Cells.Item(2, 2).Value = "LOW BAU"
$ws.Cells.Item(2, 3).Value = 2025
$ws.Cells.Item(2, 4).Value = -743.09711157562333
$ws.Cells.Item(3, 2).Value = "LOW BAU"
$ws.Cells.Item(3, 3).Value = 2030
$ws.Cells.Item(3, 4).Value = -743.09711157562333
$ws.Cells.Item(4, 2).Value = "LOW BAU"
$ws.Cells.Item(4, 3).Value = 2035
$ws.Cells.Item(4, 4).Value = -659.26322942720196
$ws.Cells.Item(5, 2).Value = "LOW BAU"
$ws.Cells.Item(5, 3).Value = 2040
$ws.Cells.Item(5, 4).Value = -554.02036161887133
$ws.Cells.Item(6, 2).Value = "LOW BAU"
$ws.Cells.Item(6, 3).Value = 2045
$ws.Cells.Item(6, 4).Value = -509.18138093476944
$ws.Cells.Item(7, 2).Value = "LOW BAU"
$ws.Cells.Item(7, 3).Value = 2050
$ws.Cells.Item(7, 4).Value = -475.75686480370655
$ws.Cells.Item(8, 2).Value = "HIGH BAU"
$ws.Cells.Item(8, 3).Value = 2025
$ws.Cells.Item(8, 4).Value = -743.09711157562333
$ws.Cells.Item(9, 2).Value = "HIGH BAU"
$ws.Cells.Item(9, 3).Value = 2030
$ws.Cells.Item(9, 4).Value = -883.90327395269583
$ws.Cells.Item(10, 2).Value = "HIGH BAU"
$ws.Cells.Item(10, 3).Value = 2035
$ws.Cells.Item(10, 4).Value = -952.20654790539174
$ws.Cells.Item(11, 2).Value = "HIGH BAU"
$ws.Cells.Item(11, 3).Value = 2040
$ws.Cells.Item(11, 4).Value = -989.30734376972725
$ws.Cells.Item(12, 2).Value = "HIGH BAU"
$ws.Cells.Item(12, 3).Value = 2045
$ws.Cells.Item(12, 4).Value = -1026.4081396340628
$ws.Cells.Item(13, 2).Value = "HIGH BAU"
$ws.Cells.Item(13, 3).Value = 2050
$ws.Cells.Item(13, 4).Value = -1063.5089354983984
$ws.Cells.Item(14, 2).Value = "Low Range WAM"
$ws.Cells.Item(14, 3).Value = 2025
$ws.Cells.Item(14, 4).Value = -834.16236046451206
$ws.Cells.Item(15, 2).Value = "Low Range WAM"
$ws.Cells.Item(15, 3).Value = 2030
$ws.Cells.Item(15, 4).Value = -834.16236046451206
$ws.Cells.Item(16, 2).Value = "Low Range WAM"
$ws.Cells.Item(16, 3).Value = 2035
$ws.Cells.Item(16, 4).Value = -782.87394609386831
$ws.Cells.Item(17, 2).Value = "Low Range WAM"
$ws.Cells.Item(17, 3).Value = 2040
$ws.Cells.Item(17, 4).Value = -708.83046661887136
$ws.Cells.Item(18, 2).Value = "Low Range WAM"
$ws.Cells.Item(18, 3).Value = 2045
$ws.Cells.Item(18, 4).Value = -667.73587426810309
$ws.Cells.Item(19, 2).Value = "Low Range WAM"
$ws.Cells.Item(19, 3).Value = 2050
$ws.Cells.Item(19, 4).Value = -636.39865194656397
$ws.Cells.Item(20, 2).Value = "High Range WAM"
$ws.Cells.Item(20, 3).Value = 2025
$ws.Cells.Item(20, 4).Value = -1041.8755864526959
$ws.Cells.Item(21, 2).Value = "High Range WAM"
$ws.Cells.Item(21, 3).Value = 2030
$ws.Cells.Item(21, 4).Value = -1041.8755864526959
$ws.Cells.Item(22, 2).Value = "High Range WAM"
$ws.Cells.Item(22, 3).Value = 2035
$ws.Cells.Item(22, 4).Value = -1166.566857905392
$ws.Cells.Item(23, 2).Value = "High Range WAM"
$ws.Cells.Item(23, 3).Value = 2040
$ws.Cells.Item(23, 4).Value = -1261.1472212697272
$ws.Cells.Item(23, 1).Value = "lulucf"
$ws.Cells.Item(23, 5).Value = "United States"
$ws.Cells.Item(23, 6).Value = "MMT CO2e"
$ws.Cells.Item(23, 7).Value = "LandSink"
$ws.Cells.Item(24, 2).Value = "High Range WAM"
$ws.Cells.Item(24, 3).Value = 2045
$ws.Cells.Item(24, 4).Value = -1318.0142513007293
$ws.Cells.Item(24, 1).Value = "lulucf"
$ws.Cells.Item(24, 5).Value = "United States"
$ws.Cells.Item(24, 6).Value = "MMT CO2e"
$ws.Cells.Item(24, 7).Value = "LandSink"
$ws.Cells.Item(25, 2).Value = "High Range WAM"
$ws.Cells.Item(25, 3).Value = 2050
$ws.Cells.Item(25, 4).Value = -1352.8343715698268
$ws.Cells.Item(25, 1).Value = "lulucf"
$ws.Cells.Item(25, 5).Value = "United States"
$ws.Cells.Item(25, 6).Value = "MMT CO2e"
$ws.Cells.Item(25, 7).Value = "LandSink"

# Update the sheet view state (top-left visible cell + active selection)
$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("I14").Select()
